# Commit: a new row of work-log data (BSHQ240813 / 王凯 / 分析优化 / 黑色素瘤)
# was inserted above row 16, pushing the former row-16 ("BSJF240124" / 蒋镥 …)
# down to row 17 and the former row-17 ("BS.develop" transcription-factor
# task) down to the previously-blank row 18.
#
# The underlying row/cell style ids for rows 16 & 17 are unchanged by the
# edit, and A18:D18 / I18 already carry the correct styles (s="7"/"30") -
# only their (previously empty) values need to be filled in. The one real
# formatting change is E18:H18, which must pick up the date number-format
# that E17:H17 already use (style id 34) instead of the old placeholder
# style (id 20). We copy that formatting over with PasteSpecial so the
# workbook reuses the existing style record instead of minting a new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bring E18:H18's formatting in line with E17:H17 (date-formatted cells)
#    before writing values into them.
$ws.Range("E17:H17").Copy()
$ws.Range("E18:H18").PasteSpecial(-4122)   # xlPasteFormats

# 2) Shift the old row-17 content ("BS.develop" / transcription-factor task)
#    down into row 18.
$ws.Range("A18").Value = "BS.develop"
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = "模块开发"
$ws.Range("D18").Value = "转录因子数据获取和可视化"
$ws.Range("E18").Value = 45742
$ws.Range("F18").Value = 45743
$ws.Range("G18").Value = 45743
$ws.Range("H18").Value = 45743
$ws.Range("I18").Value = "抓取 hTFtarget 数据库程序，转录因子可视化程序"

# 3) Shift the old row-16 content ("BSJF240124" / 蒋镥 …) down into row 17,
#    and clear the note that used to live at I17 (it moved to I18 above).
$ws.Range("A17").Value = "BSJF240124"
$ws.Range("B17").Value = "蒋镥"
$ws.Range("C17").Value = "生信协助"
$ws.Range("D17").Value = "结合转录因子"
$ws.Range("E17").Value = 45534
$ws.Range("F17").Value = 45740
$ws.Range("G17").Value = 45743
$ws.Range("H17").Value = 45743
$ws.Range("I17").Value = ""

# 4) Write the brand-new row-16 record.
$ws.Range("A16").Value = "BSHQ240813"
$ws.Range("B16").Value = "王凯"
$ws.Range("C16").Value = "分析优化"
$ws.Range("D16").Value = "黑色素瘤"
$ws.Range("E16").Value = 45743
$ws.Range("F16").Value = 45744
$ws.Range("G16").Value = 45743
$ws.Range("H16").Value = 45743
$ws.Range("I16").Value = ""
